$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume number + date range) ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Data table updates ---
$ws.Range("N15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N15").Value = -100
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 4
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 2
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E16").Value = 100
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("F16").Value = 17
$ws.Range("G16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 12
$ws.Range("H16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H16").Value = 41.666666666666
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("I16").Value = 10
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("J16").Value = 5
$ws.Range("K16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K16").Value = 100
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M16").Value = -28.571428571428
$ws.Range("N16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N16").Value = -70.588235294117
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("C17").Value = 6
$ws.Range("E17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E17").Value = 200
$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("F17").Value = 18
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 11
$ws.Range("H17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H17").Value = 63.636363636363
$ws.Range("I17").NumberFormat = "#,##0"
$ws.Range("I17").Value = 12
$ws.Range("J17").NumberFormat = "#,##0"
$ws.Range("J17").Value = 4
$ws.Range("K17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K17").Value = 200
$ws.Range("L17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L17").Value = 9.090909090909
$ws.Range("M17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M17").Value = 20
$ws.Range("N17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N17").Value = -45.454545454545
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 1
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = -50
$ws.Range("F18").NumberFormat = "#,##0"
$ws.Range("F18").Value = 5
$ws.Range("H18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H18").Value = -68.75
$ws.Range("I18").NumberFormat = "#,##0"
$ws.Range("I18").Value = 3
$ws.Range("J18").NumberFormat = "#,##0"
$ws.Range("J18").Value = 7
$ws.Range("K18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K18").Value = -57.142857142857
$ws.Range("L18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L18").Value = -40
$ws.Range("M18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M18").Value = -25
$ws.Range("N18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N18").Value = -85.714285714285
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("C19").Value = 11
$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("D19").Value = 7
$ws.Range("E19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").NumberFormat = "#,##0"
$ws.Range("F19").Value = 37
$ws.Range("G19").NumberFormat = "#,##0"
$ws.Range("G19").Value = 38
$ws.Range("H19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H19").Value = -2.631578947368
$ws.Range("I19").NumberFormat = "#,##0"
$ws.Range("I19").Value = 21
$ws.Range("J19").NumberFormat = "#,##0"
$ws.Range("J19").Value = 13
$ws.Range("K19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K19").Value = 61.538461538461
$ws.Range("L19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L19").Value = 16.666666666666
$ws.Range("M19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M19").Value = 75
$ws.Range("N19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N19").Value = 61.538461538461
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 1
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("F20").Value = 2
$ws.Range("H20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H20").Value = -50
$ws.Range("I20").NumberFormat = "#,##0"
$ws.Range("I20").Value = 1
$ws.Range("K20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K20").Value = 0
$ws.Range("L20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L20").Value = 0
$ws.Range("M20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M20").Value = 0
$ws.Range("N20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N20").Value = -83.333333333333
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("C21").Value = 23
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 13
$ws.Range("E21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("E21").Value = 76.923076923076
$ws.Range("F21").NumberFormat = "#,##0"
$ws.Range("F21").Value = 79
$ws.Range("G21").NumberFormat = "#,##0"
$ws.Range("G21").Value = 81
$ws.Range("H21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("H21").Value = -2.469135802469
$ws.Range("I21").NumberFormat = "#,##0"
$ws.Range("I21").Value = 47
$ws.Range("J21").NumberFormat = "#,##0"
$ws.Range("J21").Value = 30
$ws.Range("K21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("K21").Value = 56.666666666666
$ws.Range("L21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("L21").Value = 14.634146341463
$ws.Range("M21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("M21").Value = 14.634146341463
$ws.Range("N21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("N21").Value = -52.040816326530
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 2
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H22").Value = 0
$ws.Range("I22").NumberFormat = "#,##0"
$ws.Range("I22").Value = 2
$ws.Range("K22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K22").Value = 100
$ws.Range("L22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L22").Value = 100
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 5
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").NumberFormat = "#,##0"
$ws.Range("F23").Value = 11
$ws.Range("H23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H23").Value = 120
$ws.Range("I23").NumberFormat = "#,##0"
$ws.Range("I23").Value = 9
$ws.Range("K23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K23").Value = 800
$ws.Range("L23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L23").Value = 200
$ws.Range("M23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M23").Value = 800
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("C24").Value = 15
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("D24").Value = 23
$ws.Range("E24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E24").Value = -34.782608695652
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("F24").Value = 60
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("G24").Value = 95
$ws.Range("H24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H24").Value = -36.842105263157
$ws.Range("I24").NumberFormat = "#,##0"
$ws.Range("I24").Value = 24
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("J24").Value = 49
$ws.Range("K24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K24").Value = -51.020408163265
$ws.Range("L24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L24").Value = -41.463414634146
$ws.Range("M24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M24").Value = -17.241379310344
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("C25").Value = 10
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 5
$ws.Range("E25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E25").Value = 100
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("F25").Value = 31
$ws.Range("G25").NumberFormat = "#,##0"
$ws.Range("G25").Value = 16
$ws.Range("H25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H25").Value = 93.75
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("I25").Value = 17
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("J25").Value = 9
$ws.Range("K25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K25").Value = 88.888888888888
$ws.Range("L25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L25").Value = -10.526315789473
$ws.Range("M25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M25").Value = -5.555555555555
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 5
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 3
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 1
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L27").Value = -83.333333333333
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("G28").Value = 1
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("G29").Value = 1
$ws.Range("J39").NumberFormat = "#,##0"
$ws.Range("J39").Value = 234
$ws.Range("K39").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K39").Value = -17.021276595744
$ws.Range("L39").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L39").Value = -46.697038724373
$ws.Range("M39").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M39").Value = -67.679558011049
$ws.Range("N39").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N39").Value = -73.883928571428
$ws.Range("J41").NumberFormat = "#,##0"
$ws.Range("J41").Value = 383
$ws.Range("K41").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K41").Value = 147.096774193548
$ws.Range("L41").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L41").Value = 54.435483870967
$ws.Range("M41").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M41").Value = -20.703933747412
$ws.Range("N41").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N41").Value = -26.487523992322
$ws.Range("J43").NumberFormat = "#,##0"
$ws.Range("J43").Value = 979
$ws.Range("K43").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K43").Value = 1.767151767151
$ws.Range("L43").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L43").Value = -34.689793195463
$ws.Range("M43").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M43").Value = -69.463505926388
$ws.Range("N43").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N43").Value = -73.265974877116
